$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44362
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 6500
$ws.Range("M2").Value = 6500
$ws.Range("P2").Value = 181
$ws.Range("D3").Value = 44372
$ws.Range("J3").Value = 150
$ws.Range("N3").Value = '$/caja 36 atados'
$ws.Range("P3").Value = 194
$ws.Range("Q3").Value = 36
$ws.Range("D4").Value = 44342
$ws.Range("D5").Value = 44369
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("N5").Value = '$/caja 20 docenas'
$ws.Range("P5").Value = 7000
$ws.Range("Q5").Value = 1
$ws.Range("D6").Value = 44355
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("P6").Value = 194
$ws.Range("D7").Value = 44358
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("D8").Value = 44354
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("O8").Value = 'Región del Maule'
$ws.Range("P8").Value = 194
$ws.Range("D9").Value = 44348
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("D10").Value = 44376
$ws.Range("J10").Value = 150
$ws.Range("D11").Value = 44364
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("N11").Value = '$/caja 36 atados'
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("P11").Value = 194
$ws.Range("Q11").Value = 36
$ws.Range("D12").Value = 44386
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 6500
$ws.Range("M12").Value = 6500
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 181
$ws.Range("D13").Value = 44340
$ws.Range("O13").Value = 'Región del Maule'
$ws.Range("D14").Value = 44357
$ws.Range("K14").Value = 6500
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6500
$ws.Range("N14").Value = '$/caja 20 docenas'
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("P14").Value = 6500
$ws.Range("Q14").Value = 1
$ws.Range("D15").Value = 44371
$ws.Range("K15").Value = 6500
$ws.Range("L15").Value = 6500
$ws.Range("M15").Value = 6500
$ws.Range("O15").Value = 'Región Metropolitana'
$ws.Range("P15").Value = 181
